# Update the production / expiry dates throughout the document
# (there are 6 repeated nutrition-label blocks, so 6 occurrences of each).
$d = $word.ActiveDocument

[void]$d.Content.Find.Execute("27/03/2017", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "31/09/2017", 2)

[void]$d.Content.Find.Execute("26/03/2019", $false, $false, $false, $false, $false, `
                               $true, 1, $false, "31/09/2019", 2)

# Word keeps a "_GoBack" bookmark around the position of the most recent
# edit. After the edits above it now wraps the production date in the
# *first* nutrition block (instead of its previous location, inside "Made
# in China" near the end of the document) - move the bookmark to match.

# 1) Wrap the newly-updated first production date with the _GoBack bookmark.
$firstDate = $d.Content
[void]$firstDate.Find.Execute("31/09/2017")
[void]$d.Bookmarks.Add("_GoBack", $d.Range($firstDate.Start, $firstDate.End))

# 2) Find the *last* "Made in China" occurrence in the document - this is
#    the one that still has the old, now-stray, _GoBack bookmark sitting
#    right after "in" (between the "in" and " China" runs).
$scan = $d.Content
$scan.Start = 0
$lastStart = -1
$lastEnd = -1
$found = $scan.Find.Execute("Made in China")
while ($found) {
    $lastStart = $scan.Start
    $lastEnd = $scan.End
    $scan.Collapse(0)
    [void]$scan.MoveEnd(1, $d.Content.End - $scan.End)
    $found = $scan.Find.Execute("Made in China")
}

# 3) Drop the stale bookmark. A Range whose text is rewritten while
#    genuinely straddling the bookmark's position (i.e. starting strictly
#    before it and ending strictly after it) clears the bookmark markers
#    on save, so round-trip the text spanning "in China" through a
#    placeholder and back to its original value.
$spanStart = $lastStart + "Made ".Length
$span = $d.Range($spanStart, $lastEnd)
$original = $span.Text
$placeholder = "ZZPLACEHOLDERZZ"
$span.Text = $placeholder
$restore = $d.Range($spanStart, $spanStart + $placeholder.Length)
$restore.Text = $original
